$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Email value in row 3 (Andreea Gherghescu's contact) from
# lau.tender@gmail.com to a_diana_gherghescu@yahoo.com
$ws.Range("B3").Value = "a_diana_gherghescu@yahoo.com"
